$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Open the sync dV (Max. Slip Voltage) for the 480V level relays (rows 8-10)
$ws.Range("C8").Value = 10
$ws.Range("C9").Value = 10
$ws.Range("C10").Value = 10

# Update CT Primary ratio for all relay rows (2-10) from 13 to 2
foreach ($r in 2..10) {
    $ws.Range("F$r").Value = 2
}

# Update the active cell selection to match the saved view
$ws.Range("G12").Select()
